$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update D/E columns (start/end year) per diff
$ws.Range("D2").Value = 2010
$ws.Range("E2").Value = 2025

$ws.Range("D3").Value = 2015
$ws.Range("E3").Value = 2025

$ws.Range("D4").Value = 2019
$ws.Range("E4").Value = 2025

$ws.Range("D5").Value = 2024
$ws.Range("E5").Value = 2025

$ws.Range("D6").Value = 2010
$ws.Range("E6").Value = 2025

$ws.Range("D7").Value = 2016
$ws.Range("E7").Value = 2025

$ws.Range("D8").Value = 2019
$ws.Range("E8").Value = 2025

$ws.Range("D9").Value = 2024
$ws.Range("E9").Value = 2025

$ws.Range("D10").Value = 2024
$ws.Range("E10").Value = 2025

# Update the active cell selection to E15
$ws.Range("E15").Select()
